$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data (was Indo Farm Tractors / 1020 DI), now Sonalika Tractors / Tiger DI 50
$ws.Range("A2").Value = "Sonalika Tractors"
$ws.Range("B2").Value = "Tiger DI 50"
$ws.Range("C2").Value = "['TigerDI50img0-sonalika-tiger-di-50-1696592364.png', 'TigerDI50img1-sonalika-tiger-di-50-16965923640.png', 'TigerDI50img2-sonalika-tiger-di-50-1696592364.png']"

# Add a new row 3 for the second tractor model (Tiger 47)
$ws.Range("A3").Value = "Sonalika Tractors"
$ws.Range("B3").Value = "Tiger 47"
$ws.Range("C3").Value = "['Tiger47img0-tiger-47-1631530246.png', 'Tiger47img1-upload-1631530246-0.png', 'Tiger47img2-tiger-47-1631530246.png']"
